$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'crop athletic leggings',
    'kid basketball knee pads',
    'men''s spandex leggings',
    'hex pants',
    'medical compression pants',
    'fitness tights for men',
    'the rock mens basketball',
    'black men tights',
    'd man basketball',
    'asics compression pants',
    'knee armor knee pads',
    'adidas tights for men',
    'slide on knee pads',
    'men athletic tights',
    'blue mens compression pants',
    'sport tights for men',
    'compression basketball tights',
    'kids compression pants',
    'usa tights men',
    'razor knee pads',
    'knee pads addidas',
    'knee pads hunting',
    'knee pads leggings',
    'knee pads elbow pads youth',
    'knee pads slim',
    'basketball pants adidas',
    'mens tights navy',
    'mens adidas basketball pants',
    'compression pants men 3xl',
    'compression pants 2xu',
    'pro x knee pad',
    'puma compression pants men',
    'men''s basketball pants',
    'gray baseball pants youth girls',
    'mens basketball jacket',
    'photography knee pads',
    'men''s tights leggings',
    'baseball pants men grey',
    'venom compression pants',
    'padded knee tights',
    'ua compression pants',
    'men workout tights',
    'flag compression pants',
    'cool knee pads',
    'navy compression leggings',
    'force knee pads',
    'mens wrestling pants',
    'mens capri compression pants',
    'woman compression pants',
    'purple knee pads',
    'reebok knee pads',
    'venum compression pants men',
    'purple athletic leggings',
    'thermal compression pants',
    'addidas knee pads',
    'jordan mens tights',
    'purple compression pants men',
    'russell compression pants',
    'blue knee pads for basketball',
    'elbow knee pad',
    'tights mens',
    'super compression leggings',
    'mens leggins',
    'knee pads for teens',
    'green mens compression pants',
    'pants with padded knees',
    'compression with pads',
    'knee pads for sleeping',
    'mens winter compression pants',
    'knee pads for hvac',
    'yoga pants with knee pads',
    'black pants with knee pads',
    'kids compression knee pads',
    'elite basketball pants',
    'nike leggings mens',
    'compression pants men 3 pack',
    'compression pants baseball',
    'colored compression leggings',
    'mens compression pants with pockets',
    'mens compression pants xxl',
    'knee pads wheels',
    'ua basketball knee pads',
    'protective knee pad',
    'grey knee pads basketball',
    'orange knee pads for basketball',
    'adidas youth compression pants',
    'copper compression tights for men',
    'basketball knee pads youth boys mcdavid',
    'nike youth basketball knee pads',
    'nike kneepads',
    'knee pad strap',
    'internal knee pad',
    'basketball legings',
    'basketball knee sleves',
    'under armour baseball pants men',
    'cold gear compression leggings men',
    'mcgregor knee pads',
    'knee compression pants men',
    'compression knee pads for basketball for kids',
    'tesla compression pant',
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}